$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

# Remove the last data row (old row 22); the whole table will be
# rewritten below, so delete it first so the dimension collapses to
# A1:E21 once the remaining rows are overwritten with new content.
$ws1.Rows.Item(22).Delete()

$ws1.Range("A2").Value = "Última actualización: 05:57:04"
$ws1.Range("A3").Value = "Total filas: 16"

$rows1 = @(
    @("05:57:04","06:09","10_OLMOS",12,"LP1912"),
    @("05:57:04","06:16","215A_EL PATO",19,"LP1912"),
    @("05:57:04","06:30","23_HERNANDEZ",33,"LP1912"),
    @("05:57:04","06:34","11_ETCHEVERRY",37,"LP1912"),
    @("05:57:04","06:39","17X38_ROMERO",42,"LP1912"),
    @("05:57:04","06:41","16_SANTA ANA",44,"LP1912"),
    @("05:57:04","06:57","215A_EL PATO",60,"LP1912"),
    @("05:57:04","06:59","225_GOMEZ",62,"LP1912"),
    @("05:57:04","07:16","215C_EL PATO",79,"LP1912"),
    @("05:57:04","07:19","14_ABASTO",82,"LP1912"),
    @("05:57:04","07:21","16_SANTA ANA",84,"LP1912"),
    @("05:57:04","07:22","23_HERNANDEZ",85,"LP1912"),
    @("05:57:04","07:29","17X38_ROMERO",92,"LP1912"),
    @("05:57:04","07:35","10_OLMOS",98,"LP1912"),
    @("05:57:04","07:37","27_EL RETIRO",100,"LP1912"),
    @("05:57:04","07:55","14_ABASTO",118,"LP1912")
)

$r = 6
foreach ($row in $rows1) {
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# Sheet 2: LP1912-215
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 05:57:04"

$rows2 = @(
    @("05:57:04","06:16","215A_EL PATO",19,"LP1912"),
    @("05:57:04","06:57","215A_EL PATO",60,"LP1912"),
    @("05:57:04","07:16","215C_EL PATO",79,"LP1912")
)

$r = 6
foreach ($row in $rows2) {
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
    $ws2.Cells.Item($r, 4).Value = $row[3]
    $ws2.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# Sheet 3: 6203-6173
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 05:57:04"
$ws3.Range("A3").Value = "Total filas: 1"

$ws3.Range("A5").Value = "Hora_Scrap"
$ws3.Range("B5").Value = "Hora_Llegada"
$ws3.Range("C5").Value = "Linea"
$ws3.Range("D5").Value = "Minutos"
$ws3.Range("E5").Value = "Parada"
$ws1.Range("A5:E5").Copy()
$ws3.Range("A5:E5").PasteSpecial(-4122)

$ws3.Range("A6").Value = "05:57:04"
$ws3.Range("B6").Value = "07:43"
$ws3.Range("C6").Value = "215A_LA PLATA"
$ws3.Range("D6").Value = 106
$ws3.Range("E6").Value = "L6173"
